$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit: "137, my idea of mosque"
# A new blog entry (ser: 137) is inserted into the Day-9 feed, shifting
# the existing blog entries that were in C9 and G9 down into G9 and I9.

$oldC9 = $ws.Range("C9").Value()
$oldG9 = $ws.Range("G9").Value()

$newBlog137 = "type: blog" + "`n" + "width: 2" + "`n" + "height: 1" + "`n" + "ser: 137"

$ws.Range("I9").Value = $oldG9
$ws.Range("G9").Value = $oldC9
$ws.Range("C9").Value = $newBlog137

# Best-effort: scroll the view so column D is the left-most visible column
# (matches topLeftCell moving from E8 to D8 in the saved sheetView).
$win = $wb.Windows.Item(1)
$win.ScrollColumn = 4
$win.ScrollRow = 8
